$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-12: update Price (D) and/or Volume(1h) (E) only
$ws.Range("D2").Value = "68.633.86"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "2.464.94"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "559.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.507"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -0.51%  "
$ws.Range("E10").Value = "  +0.59%  "
$ws.Range("E11").Value = "  -2.64%  "
$ws.Range("E12").Value = "  +0.86%  "

# Rows 13-51: full row refresh (Coin, Link, Price, Volume(1h))
$ws.Range("B13").Value = "WrappedBTC"
$ws.Range("C13").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D13").Value = "68.575.89"
$ws.Range("E13").Value = "  +0.47%  "
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000168"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.82%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.33"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.08%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "10.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.90%  "
$ws.Range("B17").Value = "BitcoinCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "334.55"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.61%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.33%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.99%  "
$ws.Range("B20").Value = "SuiNetwork"
$ws.Range("C20").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.21%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("B22").Value = "Litecoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.41%  "
$ws.Range("B23").Value = "NEARProtocol"
$ws.Range("C23").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.64"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.91%  "
$ws.Range("B24").Value = "Aptos"
$ws.Range("C24").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.92%  "
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").Value = "0.0₃0817"
$ws.Range("E25").Value = "  -2.73%  "
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.47%  "
$ws.Range("B27").Value = "FirstDigitalUSD"
$ws.Range("C27").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("B28").Value = "Bittensor"
$ws.Range("C28").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "429.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.55%  "
$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.14"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.43%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.61"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.90%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.06"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.43%  "
$ws.Range("B32").Value = "WhiteBITCoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.110"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.11%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.74"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.87%  "
$ws.Range("B36").Value = "PolygonEcosystemToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.300"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.21%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.41"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.85%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.54%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.08"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.21%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.06"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.41%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.34"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.77%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "129.15"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.26%  "
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0715"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.23%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.482"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.83%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.560"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0910"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.44%  "
$ws.Range("B47").Value = "BitgetToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.14%  "
$ws.Range("B48").Value = "Optimism"
$ws.Range("C48").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.39"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.15%  "
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.49%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.78"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.90%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0206"
$ws.Range("E51").Value = "  +0.77%  "
